$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values per row (B, C, D, E, G); F unchanged
$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 13.86384647080068;  G = 21.98653043760045 }
    3 = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 3.223369029078222;  E = 13.86384647080068;  G = 20.15985084044064 }
    4 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 3.223369029078222;  E = 0.5333859586016987; G = 8.656069925401464 }
    5 = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
    6 = @{ B = 3.272327238179451;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 6.15379541431027 }
    7 = @{ B = 0.2881169905109251; C = 1.626987699542094; D = 0.1496068669990043; E = 0.5333859586016987; G = 2.598097515653722 }
    8 = @{ B = 1.445647641019636;  C = 1.626987699542094; D = 0.7210945179870265; E = 0.5333859586016987; G = 4.327115817150455 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
